# Refresh the crypto price/volume snapshot (data-only update; no rows
# added or removed). A handful of values happen to look like plain
# numbers ("579.30", "1.00", "0.0000249", ...) but the source sheet
# keeps every Price/Volume cell as text, so those are entered with a
# leading apostrophe (exactly as typing them into Excel would do) and
# then restored to the sheet's plain "Normal" look.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.959.93"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "3.101.66"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'579.30"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'172.39"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.096.96"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").Value = "'6.44"
$ws.Range("E10").Value = "  -4.53%  "
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'37.29"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "3.613.95"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").Value = "66.991.50"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "'7.19"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "3.103.14"
$ws.Range("E19").Value = "  +2.62%  "
$ws.Range("D20").Value = "'16.30"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'486.18"
$ws.Range("E21").Value = "  +3.73%  "
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.718"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'7.55"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'13.41"
$ws.Range("E24").Value = "  +4.56%  "
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'84.41"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'2.36"
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("D26").Style = "Normal"
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D27").Style = "Normal"
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.00"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'8.07"
$ws.Range("E29").Value = "  -4.93%  "
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("D31").Value = "'2.67"
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'28.95"
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'0.0000100"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'0.114"
$ws.Range("E34").Value = "  -3.29%  "
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'5.91"
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "'47.59"
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'2.11"
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'50.14"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.315"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'8.67"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'2.80"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "2.837.60"
$ws.Range("E45").Value = "  +3.98%  "
$ws.Range("D46").Value = "'0.0361"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'384.76"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'136.10"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'24.91"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "  -0.87%  "
